$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder "Recorded By" list ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Row 3: reorder "Recorded By" list ---
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 4: reorder "Recorded By" list ---
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 5: reorder "Recorded By" list ---
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Row 6: Recorded Sessions count (Class Statistics) ---
$ws.Range("L6").Value = 13

# --- Row 7: reorder "Recorded By" list + Missing Sessions count ---
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("L7").Value = 1

# --- Row 9: Coverage % (text, keep style/number-format untouched) ---
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "44.8%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# --- Row 10: Average Attendance % (text, keep style/number-format untouched) ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "24.5%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- Row 12: reorder "Recorded By" list ---
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg"

# --- Row 15: reorder "Recorded By" list + Group Statistics row ---
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("O15").Value = 13
$ws.Range("P15").Value = 1

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "44.8%"
$ws.Range("K15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "24.5%"
$ws.Range("K15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# --- Row 27: PHARMACOLOGY session 2 now recorded (copy "Recorded" row format, then fill values) ---
$ws.Range("A26:I26").Copy()
$ws.Range("A27:I27").PasteSpecial(-4122)
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg"
$ws.Range("H27").Value = "76/251"
$ws.Range("I27").Value = "Recorded"
